$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 44: new time entry (Coding) - set D (Interruption) first so the
# dependent shared formula in E44 picks up its value on recalculation.
$ws.Range("D44").Value = 5
$ws.Range("A44").Value = 41897
$ws.Range("B44").Value = 0.98333333333333339
$ws.Range("C44").Value = 1.0249999999999999
$ws.Range("F44").Value = "Coding"

# Row 45: new time entry (Testing)
$ws.Range("D45").Value = 5
$ws.Range("A45").Value = 41897
$ws.Range("B45").Value = 0.027083333333333334
$ws.Range("C45").Value = 0.04861111111111111
$ws.Range("F45").Value = "Testing"

# Update the view: scroll back to top and change the active selection
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("F21").Select()

$excel.CalculateFull()

# Make sure the pie chart on Sheet2 (driven by the SUMIF totals above)
# picks up the refreshed totals.
$ws2 = $wb.Worksheets.Item("Sheet2")
$chart = $ws2.ChartObjects().Item(1).Chart
$chart.Refresh()
